$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "67.214.27"
Set-TextValue $ws.Range("E2") "  -1.69%  "
Set-TextValue $ws.Range("D3") "2.485.64"
Set-TextValue $ws.Range("E3") "  -2.05%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.02%  "
Set-TextValue $ws.Range("D5") "584.28"
Set-TextValue $ws.Range("E5") "  -1.79%  "
Set-TextValue $ws.Range("D6") "166.29"
Set-TextValue $ws.Range("E6") "  -6.65%  "
Set-TextValue $ws.Range("E7") "  +0.14%  "
Set-TextValue $ws.Range("D8") "0.516"
Set-TextValue $ws.Range("E8") "  -3.05%  "
Set-TextValue $ws.Range("D9") "2.485.86"
Set-TextValue $ws.Range("E9") "  -2.02%  "
Set-TextValue $ws.Range("D10") "0.136"
Set-TextValue $ws.Range("E10") "  -3.78%  "
Set-TextValue $ws.Range("E11") "  +0.21%  "
Set-TextValue $ws.Range("D12") "0.341"
Set-TextValue $ws.Range("E12") "  -2.15%  "
Set-TextValue $ws.Range("D13") "4.94"
Set-TextValue $ws.Range("E13") "  -3.91%  "
Set-TextValue $ws.Range("D14") "25.92"
Set-TextValue $ws.Range("E14") "  -4.01%  "
Set-TextValue $ws.Range("D15") "2.938.54"
Set-TextValue $ws.Range("E15") "  -2.11%  "
Set-TextValue $ws.Range("D16") "0.0000173"
Set-TextValue $ws.Range("E16") "  -3.81%  "
Set-TextValue $ws.Range("D17") "66.968.84"
Set-TextValue $ws.Range("E17") "  -1.87%  "
Set-TextValue $ws.Range("D18") "2.434.60"
Set-TextValue $ws.Range("E18") "  -4.49%  "
Set-TextValue $ws.Range("D19") "11.63"
Set-TextValue $ws.Range("E19") "  +0.56%  "
Set-TextValue $ws.Range("D20") "7.82"
Set-TextValue $ws.Range("E20") "  -2.71%  "
Set-TextValue $ws.Range("D21") "360.13"
Set-TextValue $ws.Range("E22") "  -2.36%  "
Set-TextValue $ws.Range("D23") "4.41"
Set-TextValue $ws.Range("E23") "  -6.88%  "
Set-TextValue $ws.Range("E24") "  +0.08%  "
Set-TextValue $ws.Range("D25") "70.56"
Set-TextValue $ws.Range("E25") "  -0.34%  "
Set-TextValue $ws.Range("E26") "  -5.02%  "
Set-TextValue $ws.Range("D27") "9.41"
Set-TextValue $ws.Range("E27") "  -8.22%  "
Set-TextValue $ws.Range("D28") "0.997"
Set-TextValue $ws.Range("E28") "  +0.00%  "
Set-TextValue $ws.Range("D30") "0.0₃0935"
Set-TextValue $ws.Range("E30") "  -6.47%  "
Set-TextValue $ws.Range("D31") "8.03"
Set-TextValue $ws.Range("E31") "  -3.42%  "
Set-TextValue $ws.Range("D32") "498.63"
Set-TextValue $ws.Range("E32") "  -8.19%  "
Set-TextValue $ws.Range("D33") "1.83"
Set-TextValue $ws.Range("E33") "  -2.76%  "
Set-TextValue $ws.Range("D34") "1.27"
Set-TextValue $ws.Range("E34") "  -5.45%  "
Set-TextValue $ws.Range("D35") "1.00"
Set-TextValue $ws.Range("E35") "  +0.04%  "
Set-TextValue $ws.Range("D36") "0.126"
Set-TextValue $ws.Range("E36") "  -2.67%  "
Set-TextValue $ws.Range("D37") "159.41"
Set-TextValue $ws.Range("D38") "19.07"
Set-TextValue $ws.Range("E38") "  +0.94%  "
Set-TextValue $ws.Range("D39") "1.42"
Set-TextValue $ws.Range("E39") "  -3.22%  "
Set-TextValue $ws.Range("D40") "18.57"
Set-TextValue $ws.Range("E40") "  -0.72%  "
Set-TextValue $ws.Range("D41") "1.73"
Set-TextValue $ws.Range("E41") "  -4.77%  "
Set-TextValue $ws.Range("D42") "4.94"
Set-TextValue $ws.Range("E42") "  -5.73%  "
Set-TextValue $ws.Range("D43") "0.336"
Set-TextValue $ws.Range("E43") "  -5.99%  "
Set-TextValue $ws.Range("E44") "  +0.05%  "
Set-TextValue $ws.Range("D45") "2.47"
Set-TextValue $ws.Range("E45") "  -4.46%  "
Set-TextValue $ws.Range("D46") "39.34"
Set-TextValue $ws.Range("E46") "  -1.60%  "
Set-TextValue $ws.Range("D47") "141.78"
Set-TextValue $ws.Range("E47") "  -4.21%  "
Set-TextValue $ws.Range("D48") "3.62"
Set-TextValue $ws.Range("E48") "  -3.24%  "
Set-TextValue $ws.Range("D49") "0.537"
Set-TextValue $ws.Range("E49") "  -4.65%  "
Set-TextValue $ws.Range("D50") "0.0₆0264"
Set-TextValue $ws.Range("E50") "  -5.64%  "
Set-TextValue $ws.Range("D51") "1.64"
Set-TextValue $ws.Range("E51") "  -4.17%  "
